# Auto-generated PowerShell Excel COM-interop script
# Applies numeric updates to currentAveragePrice / LevePrice / LeveProfit columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled price-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 303.7143
$ws.Range("I4").Value = 187.66667
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 187.66667
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -73.66667000000001
$ws.Range("N4").Value = -1228
$ws.Range("H32").Value = 2671.6365
$ws.Range("J32").Value = 2764.45
$ws.Range("L32").Value = 2764.45
$ws.Range("N32").Value = -3416.45
$ws.Range("H70").Value = 3313.6785
$ws.Range("I70").Value = 3332.5
$ws.Range("J70").Value = 3279.8
$ws.Range("K70").Value = 9997.5
$ws.Range("L70").Value = 9839.400000000001
$ws.Range("M70").Value = -9727.5
$ws.Range("N70").Value = -10379.4
$ws.Range("H73").Value = 3313.6785
$ws.Range("I73").Value = 3332.5
$ws.Range("J73").Value = 3279.8
$ws.Range("K73").Value = 9997.5
$ws.Range("L73").Value = 9839.400000000001
$ws.Range("M73").Value = -9061.5
$ws.Range("N73").Value = -11711.4
$ws.Range("H74").Value = 6792.9287
$ws.Range("I74").Value = 6067
$ws.Range("K74").Value = 6067
$ws.Range("M74").Value = -5131
$ws.Range("H77").Value = 6792.9287
$ws.Range("I77").Value = 6067
$ws.Range("K77").Value = 30335
$ws.Range("M77").Value = -25655
$ws.Range("H100").Value = 2876.2222
$ws.Range("I100").Value = 2412.4285
$ws.Range("K100").Value = 2412.4285
$ws.Range("M100").Value = -1871.4285
$ws.Range("H106").Value = 1922.8823
$ws.Range("I106").Value = 1692.0714
$ws.Range("K106").Value = 1692.0714
$ws.Range("M106").Value = -1061.0714
$ws.Range("H132").Value = 2512.4736
$ws.Range("I132").Value = 2695.9
$ws.Range("J132").Value = 1202.2858
$ws.Range("K132").Value = 8087.700000000001
$ws.Range("L132").Value = 3606.8574
$ws.Range("M132").Value = -5557.700000000001
$ws.Range("N132").Value = -8666.857400000001
$ws.Range("H138").Value = 22729328
$ws.Range("J138").Value = 55558820
$ws.Range("L138").Value = 166676460
$ws.Range("N138").Value = -166686740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2259.5715
$ws.Range("I2").Value = 2259.5715
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2259.5715
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2146.5715
$ws.Range("H45").Value = 1799.625
$ws.Range("I45").Value = 1000
$ws.Range("J45").Value = 4198.5
$ws.Range("K45").Value = 1000
$ws.Range("L45").Value = 4198.5
$ws.Range("M45").Value = -623
$ws.Range("N45").Value = -4952.5
$ws.Range("H116").Value = 2259.5715
$ws.Range("I116").Value = 2259.5715
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2259.5715
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 34.42849999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2259.5715
$ws.Range("I3").Value = 2259.5715
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2259.5715
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -2145.5715
$ws.Range("H22").Value = 305
$ws.Range("J22").Value = 198
$ws.Range("L22").Value = 198
$ws.Range("N22").Value = -544
$ws.Range("H86").Value = 15695.174
$ws.Range("I86").Value = 10426.417
$ws.Range("J86").Value = 21442.908
$ws.Range("K86").Value = 10426.417
$ws.Range("L86").Value = 21442.908
$ws.Range("M86").Value = -9303.416999999999
$ws.Range("N86").Value = -23688.908
$ws.Range("H89").Value = 15695.174
$ws.Range("I89").Value = 10426.417
$ws.Range("J89").Value = 21442.908
$ws.Range("K89").Value = 52132.085
$ws.Range("L89").Value = 107214.54
$ws.Range("M89").Value = -46516.085
$ws.Range("N89").Value = -118446.54
$ws.Range("H105").Value = 53007.25
$ws.Range("I105").Value = 67343
$ws.Range("K105").Value = 67343
$ws.Range("M105").Value = -65596
$ws.Range("H134").Value = 1113.92
$ws.Range("I134").Value = 1113.92
$ws.Range("K134").Value = 3341.76
$ws.Range("M134").Value = -806.7600000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4684.647
$ws.Range("I99").Value = 5157.154
$ws.Range("J99").Value = 3149
$ws.Range("K99").Value = 5157.154
$ws.Range("L99").Value = 3149
$ws.Range("M99").Value = -3659.154
$ws.Range("N99").Value = -6145
$ws.Range("H105").Value = 14660.429
$ws.Range("I105").Value = 8280.454
$ws.Range("K105").Value = 8280.454
$ws.Range("M105").Value = -6533.454
$ws.Range("H126").Value = 4684.647
$ws.Range("I126").Value = 5157.154
$ws.Range("J126").Value = 3149
$ws.Range("K126").Value = 15471.462
$ws.Range("L126").Value = 9447
$ws.Range("M126").Value = -13001.462
$ws.Range("N126").Value = -14387
$ws.Range("H131").Value = 28065.375
$ws.Range("J131").Value = 30646.143
$ws.Range("L131").Value = 30646.143
$ws.Range("N131").Value = -40726.143
$ws.Range("H132").Value = 4395.9644
$ws.Range("I132").Value = 4125.8096
$ws.Range("K132").Value = 12377.4288
$ws.Range("M132").Value = -9847.428799999998
$ws.Range("H141").Value = 85395.55
$ws.Range("J141").Value = 90995.05499999999
$ws.Range("L141").Value = 90995.05499999999
$ws.Range("N141").Value = -101355.055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2923
$ws.Range("I5").Value = 1269.1428
$ws.Range("J5").Value = 14500
$ws.Range("K5").Value = 3807.4284
$ws.Range("L5").Value = 43500
$ws.Range("M5").Value = -3695.4284
$ws.Range("N5").Value = -43724
$ws.Range("H55").Value = 6307068
$ws.Range("J55").Value = 11113943
$ws.Range("L55").Value = 33341829
$ws.Range("N55").Value = -33342183
$ws.Range("H80").Value = 25716.334
$ws.Range("I80").Value = 4000.5
$ws.Range("J80").Value = 69148
$ws.Range("K80").Value = 12001.5
$ws.Range("L80").Value = 207444
$ws.Range("M80").Value = -11065.5
$ws.Range("N80").Value = -209316
$ws.Range("H83").Value = 25716.334
$ws.Range("I83").Value = 4000.5
$ws.Range("J83").Value = 69148
$ws.Range("K83").Value = 36004.5
$ws.Range("L83").Value = 622332
$ws.Range("M83").Value = -31324.5
$ws.Range("N83").Value = -631692
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H135").Value = 2923
$ws.Range("I135").Value = 1269.1428
$ws.Range("J135").Value = 14500
$ws.Range("K135").Value = 11422.2852
$ws.Range("L135").Value = 130500
$ws.Range("M135").Value = -8887.2852
$ws.Range("N135").Value = -135570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2639
$ws.Range("J80").Value = 3131.6667
$ws.Range("L80").Value = 3131.6667
$ws.Range("N80").Value = -5127.6667
$ws.Range("H83").Value = 2639
$ws.Range("J83").Value = 3131.6667
$ws.Range("L83").Value = 15658.3335
$ws.Range("N83").Value = -25642.3335
$ws.Range("H122").Value = 1251.2
$ws.Range("I122").Value = 1034.409
$ws.Range("K122").Value = 3103.227
$ws.Range("M122").Value = -653.2270000000003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1785.3103
$ws.Range("I16").Value = 1780.1904
$ws.Range("J16").Value = 1798.75
$ws.Range("K16").Value = 1780.1904
$ws.Range("L16").Value = 1798.75
$ws.Range("M16").Value = -1610.1904
$ws.Range("N16").Value = -2138.75
$ws.Range("H22").Value = 1445.415
$ws.Range("I22").Value = 1232.1471
$ws.Range("K22").Value = 1232.1471
$ws.Range("M22").Value = -937.1470999999999
$ws.Range("H27").Value = 1445.415
$ws.Range("I27").Value = 1232.1471
$ws.Range("K27").Value = 1232.1471
$ws.Range("M27").Value = -1125.1471
$ws.Range("H53").Value = 22000
$ws.Range("I53").Value = 22000
$ws.Range("K53").Value = 22000
$ws.Range("M53").Value = -21482
$ws.Range("H68").Value = 2251.0715
$ws.Range("I68").Value = 2074.0908
$ws.Range("J68").Value = 2900
$ws.Range("K68").Value = 2074.0908
$ws.Range("L68").Value = 2900
$ws.Range("M68").Value = -1325.0908
$ws.Range("N68").Value = -4398
$ws.Range("H71").Value = 2251.0715
$ws.Range("I71").Value = 2074.0908
$ws.Range("J71").Value = 2900
$ws.Range("K71").Value = 10370.454
$ws.Range("L71").Value = 14500
$ws.Range("M71").Value = -6626.454
$ws.Range("N71").Value = -21988
$ws.Range("H132").Value = 7339.0586
$ws.Range("I132").Value = 2235.25
$ws.Range("K132").Value = 6705.75
$ws.Range("M132").Value = -4175.75
$ws.Range("H140").Value = 134000
$ws.Range("J140").Value = 130000
$ws.Range("L140").Value = 130000
$ws.Range("N140").Value = -140360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1212.2727
$ws.Range("I107").Value = 1189.0834
$ws.Range("J107").Value = 1240.1
$ws.Range("K107").Value = 3567.2502
$ws.Range("L107").Value = 3720.3
$ws.Range("M107").Value = -1647.2502
$ws.Range("N107").Value = -7560.299999999999
$ws.Range("H122").Value = 6187.3887
$ws.Range("I122").Value = 6645
$ws.Range("K122").Value = 19935
$ws.Range("M122").Value = -17485
$ws.Range("H126").Value = 6499.6665
$ws.Range("I126").Value = 4749.75
$ws.Range("J126").Value = 9999.5
$ws.Range("K126").Value = 14249.25
$ws.Range("L126").Value = 29998.5
$ws.Range("M126").Value = -11779.25
$ws.Range("N126").Value = -34938.5
